# Apply the "automatic update" edit:
#  - Column C (Förändrad) on rows 2-9: 46072 -> 46073
#  - Rows 4, 5, 9 rotate their record-specific data (A, B, F, G)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" date serial in column C for all data rows (2-9)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = 46073
}

# 2) Rotate the per-record fields among rows 4, 5 and 9
#    old row4 -> new row9, old row5 -> new row4, old row9 -> new row5
$row4 = @{
    A = $ws.Cells.Item(4, 1).Value2
    B = $ws.Cells.Item(4, 2).Value2
    F = $ws.Cells.Item(4, 6).Value2
    G = $ws.Cells.Item(4, 7).Value2
}
$row5 = @{
    A = $ws.Cells.Item(5, 1).Value2
    B = $ws.Cells.Item(5, 2).Value2
    F = $ws.Cells.Item(5, 6).Value2
    G = $ws.Cells.Item(5, 7).Value2
}
$row9 = @{
    A = $ws.Cells.Item(9, 1).Value2
    B = $ws.Cells.Item(9, 2).Value2
    F = $ws.Cells.Item(9, 6).Value2
    G = $ws.Cells.Item(9, 7).Value2
}

# New row 4 <= old row 5
$ws.Cells.Item(4, 1).Value = $row5.A
$ws.Cells.Item(4, 2).Value = $row5.B
$ws.Cells.Item(4, 6).Value = $row5.F
$ws.Cells.Item(4, 7).Value = $row5.G

# New row 5 <= old row 9
$ws.Cells.Item(5, 1).Value = $row9.A
$ws.Cells.Item(5, 2).Value = $row9.B
$ws.Cells.Item(5, 6).Value = $null
$ws.Cells.Item(5, 7).Value = $row9.G

# New row 9 <= old row 4
$ws.Cells.Item(9, 1).Value = $row4.A
$ws.Cells.Item(9, 2).Value = $row4.B
$ws.Cells.Item(9, 7).Value = $row4.G
